# Weekly update: insert a new daily record for "Cebollín" at row 17
# (Vega Monumental Concepción market), pushing all subsequent rows
# down by one. Dimension grows from A1:R50 to A1:R51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 17, shifting rows 17-50 -> 18-51.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new record's data.
$ws.Cells.Item(17, 1).Value  = 11
$ws.Cells.Item(17, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(17, 3).Value  = "Bíobío"
$ws.Cells.Item(17, 4).Value  = "2022-06-08"
$ws.Cells.Item(17, 5).Value  = 8
$ws.Cells.Item(17, 6).Value  = 100112037
$ws.Cells.Item(17, 7).Value  = "Cebollín"
$ws.Cells.Item(17, 8).Value  = "Sin especificar"
$ws.Cells.Item(17, 9).Value  = "Primera"
$ws.Cells.Item(17, 10).Value = 160
$ws.Cells.Item(17, 11).Value = 6000
$ws.Cells.Item(17, 12).Value = 6500
$ws.Cells.Item(17, 13).Value = 6250
$ws.Cells.Item(17, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(17, 15).Value = "Región Metropolitana"
$ws.Cells.Item(17, 16).Value = 174
$ws.Cells.Item(17, 17).Value = 36
$ws.Cells.Item(17, 18).Value = "Hortaliza"
